$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Expected "

$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = RGB(255, 255, 0)
